# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) for rows 2-16, replacing the previous Strike# derived values
$kValues = @{
    2  = 3
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 2
    14 = 3
    15 = 0
    16 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
